# Generate Report for Handoff
#
# A fresh handoff run completed for the row whose source file is
# "45bee6c3-3507-4a91-8d6d-41e47c2e86da.md". That refreshes two timestamp
# strings that the report happens to re-use (via the shared-string table)
# across several other rows that were stamped by the very same handoff
# batch, and it marks all of those rows' Priority column as "ht" on both
# language sheets.
#
# Rows affected on zh-cn / de-de: 8, 9, 10, 12, 13, 14 (row 11 belongs to a
# different handoff batch and is left untouched).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$rows = @(8, 9, 10, 12, 13, 14)

foreach ($r in $rows) {
    # Latest HO Xliff Generate Date (Overview) / Latest Handback DateTime (de-de)
    # both displayed "2016-08-31 12:24:51" for this handoff batch; bump to the
    # new generation time.
    $overview.Range("G$r").Value = "2016-08-31 12:25:23"
    $dede.Range("H$r").Value     = "2016-08-31 12:25:23"

    # Latest Handoff Datetime on zh-cn displayed "2016-08-31 12:24:44" for this
    # batch; bump to the new handoff time.
    $zhcn.Range("H$r").Value = "2016-08-31 12:25:18"

    # Priority column: mark these rows as belonging to the "ht" handoff type.
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"
}
